$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 with revised data
$data = @(
    @(1,  "Sandwich",         "Veg",     250, 2,  3, "accepted"),
    @(2,  "Bread",            "Veg",     100, 1,  1, "accepted"),
    @(3,  "Pizza",            "Veg",     500, 1,  2, "accepted"),
    @(4,  "Beef Burger",      "Non-Veg", 250, 2,  3, "accepted"),
    @(5,  "Chciken Pizza",    "Non-Veg", 150, 1,  1, "accepted"),
    @(6,  "Salad",            "Veg",     200, 1,  1, "accepted"),
    @(7,  "Chicken Burger",   "Non-Veg", 200, 2,  2, "accepted"),
    @(8,  "Fries",            "Veg",     500, 1,  2, "accepted"),
    @(9,  "Chicken Burrito",  "Non-Veg", 500, 2,  1, "pending"),
    @(10, "Chicken Noodeles", "Non-Veg", 210, 10, 5, "pending"),
    @(11, "KungPao Chicken",  "Non-Veg", 100, 5,  3, "accepted"),
    @(12, "Mashed Potato",    "Veg",     120, 3,  6, "pending")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
